$d = $word.ActiveDocument

# Locate the paragraph that holds the "LOB1012: Estatística (Requisito fraco)"
# text so the edit is anchored to content rather than a hard-coded index.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*LOB1012: Estatística (Requisito fraco)*") {
        $anchor = $p
        break
    }
}

# The three paragraphs immediately following the anchor are:
#   1) an empty "Normal" paragraph
#   2) an empty "Normal" paragraph with PageBreakBefore (and centered/left jc)
#   3) the "© 2020 . Contact: ..." copyright paragraph
# All three are removed, leaving the anchor paragraph directly followed by
# the remaining empty paragraph + trailing page-break paragraph.
$first = $anchor.Next()
$last = $first.Next().Next()

$range = $d.Range($first.Range.Start, $last.Range.End)
$range.Delete()
